$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Row 9 = Date property -> update the date value (column B)
$ws.Range("B9").Value = "2026-01-15T08:54:26+00:00"

# Row 12 = Jurisdiction property -> set value to FRANCE (column B)
$ws.Range("B12").Value = "FRANCE"
